$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of credential data: invalid email / invalid password / "no" (not valid)
$ws.Range("A3").Value = "invalid@test.com"
$ws.Range("B3").Value = "invalid!23"
$ws.Range("C3").Value = "no"

# Hyperlink the new email address cell, same way A2 already is
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:invalid@test.com")
$ws.Range("A3").Style = $ws.Range("A2").Style

# Column A needs to grow a bit to fit the longer e-mail address
$ws.Columns.Item(1).ColumnWidth = 15.877604166666666

$ws.Range("C3").Select()
